# Project Analysis Presentation — milestone2 edit
#   1) Slide 2 ("Project Management"): merge the two runs of the first
#      bullet ("Current iteration #2 will end on " + "Friday") into a
#      single run, and drop the stray trailing <a:endParaRPr/> that sat
#      after the "Detailed UC defined, ..." bullet.
#   2) Slide 4 ("Domain Model"): enlarge/reposition the domain-model
#      diagram picture.

$p = $ppt.ActivePresentation

# --- 1) Slide 2: tidy up the bullet text runs --------------------------
$slide2 = $p.Slides.Item(2)
$body   = $slide2.Shapes.Item(2)
$tr     = $body.TextFrame.TextRange

# Fix the "Detailed UC ..." paragraph (#2) first: delete it (text + its
# paragraph mark, which carries the unwanted trailing endParaRPr) and
# retype it, restoring its level-1 indent.
$para2 = $tr.Paragraphs(2)
$para2.Delete()
$para2 = $tr.Paragraphs(2)
$para2.InsertBefore("Detailed UC defined, Domain Model created and Project Analysis`r")
$para2.IndentLevel = 2

# Now fix the "Current iteration ..." paragraph (#1): delete it and
# retype the combined text as one run.
$para1 = $tr.Paragraphs(1)
$para1.Delete()
$para1 = $tr.Paragraphs(1)
$para1.InsertBefore("Current iteration #2 will end on Friday`r")

# Deleting paragraph 1 can reset paragraph 2's indent level again, so
# restore it once more just in case.
$para2 = $tr.Paragraphs(2)
$para2.IndentLevel = 2

# --- 2) Slide 4: resize/reposition the domain model picture ------------
$slide4 = $p.Slides.Item(4)
$pic    = $slide4.Shapes.Item(2)

$pic.Left   = 388.3471553543307
$pic.Top    = 0
$pic.Width  = 579.9284401968504
$pic.Height = 539.99997
